$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 / Location: prepend an intro clause to the address response in C2.
$ws.Range("C2").Value = "The University is located at, Manipal University Jaipur, Dehmi Kalan, Off Jaipur-Ajmer Expressway, Jaipur, (Raj.) Rajasthan 303007. Phone: +91 141-3999100 More details can be found on: https://jaipur.manipal.edu/muj/contact-us.html"

# Row 4 / Timings: append additional paraphrased questions to the Questions list in B4.
$ws.Range("B4").Value = "What is the timing of the College?, What time the college is open?, May I know the time I can contact Manipal University Jaipur?, How long the college will stay open?, What are the timings of MUJ?, At what time is the college open?, What is the time I can contact the University?"

# Row 3 / Owner: append additional paraphrased questions to the Questions list in B3.
$ws.Range("B3").Value = "Who is the President of this College?, Who is the head of this University?, What is the name of the President of Manipal University Jaipur?, Who is the Dead of this College?, Does this college have a head?, Who is the head of this College?, Who are the managers of Manipal University Jaipur?, Who is the President of MUJ?, Who is the president of this University?, Is there a board of directors in this college?"

# Row heights grow to fit the longer wrapped text (values taken from the authored workbook).
$ws.Rows.Item(2).RowHeight = 115.2
$ws.Rows.Item(3).RowHeight = 144
$ws.Rows.Item(4).RowHeight = 100.8

# Selection moved to B3 after the edit.
$ws.Range("B3").Select() | Out-Null
